# Applies the diff: fills in the missing B313 symbol cell, and appends four
# new rows (314-317) that mirror the existing data pattern for this series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 313 was missing its "symbol" (column B) value; fill it in ---
$ws.Range("B313").Value = "ECONOMICS:CNCBBS"

# --- Copy column-A's date formatting (style s="2") down onto the new rows ---
$ws.Range("A313").Copy()
$ws.Range("A314:A317").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- New rows 314-316: same pattern as row 313, including the symbol ---
for ($r = 314; $r -le 316; $r++) {
    $ws.Cells.Item($r, 1).Value = 45230
    $ws.Cells.Item($r, 2).Value = "ECONOMICS:CNCBBS"
    $ws.Cells.Item($r, 3).Value = 43325980000000
    $ws.Cells.Item($r, 4).Value = 43325980000000
    $ws.Cells.Item($r, 5).Value = 43325980000000
    $ws.Cells.Item($r, 6).Value = 43325980000000
    $ws.Cells.Item($r, 7).Value = 0
}

# --- New row 317: same pattern, but (like the original row 313) no symbol ---
$ws.Cells.Item(317, 1).Value = 45230
$ws.Cells.Item(317, 3).Value = 43325980000000
$ws.Cells.Item(317, 4).Value = 43325980000000
$ws.Cells.Item(317, 5).Value = 43325980000000
$ws.Cells.Item(317, 6).Value = 43325980000000
$ws.Cells.Item(317, 7).Value = 0
